$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") rows 2-54 currently hold the date serial 46081
# (2026-02-28). Bump each of those cells forward one day to serial 46082
# (2026-03-01), matching the upstream "Automatic update" refresh.
for ($row = 2; $row -le 54; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 46081) {
        $cell.Value2 = 46082
    }
}
